$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (B3): 1.0.0 -> 1.1.0
$ws.Range("B3").Value = "1.1.0"

# Update Date value (B8): 2023-06-07T11:52:14+02:00 -> 2023-07-10T23:08:03+02:00
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
